$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Each multiplication problem lives in its own table cell, so we scope every
# Find/Replace to that cell's Range. This keeps the operation unambiguous even
# though one new value ("36×56=") duplicates a value that existed elsewhere
# in the table before this edit - a plain document-wide replace-all could
# otherwise clobber the wrong cell depending on execution order.
$replacements = @(
    @{ Row = 1;  Col = 1; Old = "36×99="; New = "77×17=" },
    @{ Row = 1;  Col = 2; Old = "61×97="; New = "92×84=" },
    @{ Row = 1;  Col = 3; Old = "82×31="; New = "50×49=" },
    @{ Row = 1;  Col = 4; Old = "13×94="; New = "36×21=" },
    @{ Row = 1;  Col = 5; Old = "50×68="; New = "47×54=" },
    @{ Row = 5;  Col = 1; Old = "71×85="; New = "39×26=" },
    @{ Row = 5;  Col = 2; Old = "57×98="; New = "48×78=" },
    @{ Row = 5;  Col = 3; Old = "87×51="; New = "68×75=" },
    @{ Row = 5;  Col = 4; Old = "75×84="; New = "47×70=" },
    @{ Row = 5;  Col = 5; Old = "74×26="; New = "92×64=" },
    @{ Row = 10; Col = 1; Old = "44×79="; New = "58×68=" },
    @{ Row = 10; Col = 2; Old = "23×64="; New = "31×76=" },
    @{ Row = 10; Col = 3; Old = "12×81="; New = "51×37=" },
    @{ Row = 10; Col = 4; Old = "44×35="; New = "82×85=" },
    @{ Row = 10; Col = 5; Old = "50×75="; New = "22×26=" },
    @{ Row = 15; Col = 1; Old = "21×20="; New = "94×42=" },
    @{ Row = 15; Col = 2; Old = "59×72="; New = "54×95=" },
    @{ Row = 15; Col = 3; Old = "72×96="; New = "65×31=" },
    @{ Row = 15; Col = 4; Old = "56×97="; New = "60×81=" },
    @{ Row = 15; Col = 5; Old = "63×40="; New = "28×54=" },
    @{ Row = 20; Col = 1; Old = "36×56="; New = "93×34=" },
    @{ Row = 20; Col = 2; Old = "94×48="; New = "88×32=" },
    @{ Row = 20; Col = 3; Old = "21×92="; New = "78×69=" },
    @{ Row = 20; Col = 4; Old = "75×35="; New = "35×54=" },
    @{ Row = 20; Col = 5; Old = "82×88="; New = "36×56=" }
)

foreach ($item in $replacements) {
    $rng = $t.Cell($item.Row, $item.Col).Range
    $found = $rng.Find.Execute($item.Old, $false, $false, $false, $false, $false, $true, 1, $false, $item.New, 2)
    if (-not $found) {
        Write-Host "WARNING: replacement not found for" $item.Row $item.Col $item.Old
    }
}
